{"js": "// The document contains a merge-field placeholder \"{tempatLahir}, {ttl}\"\n// that must be fixed to \"{tempatLahirPembuat}, {ttlPembuat}\" so the field\n// names match the \"...Pembuat\" convention used elsewhere in the document\n// (e.g. {namaPembuat}, {alamatPembuat}, {pekerjaanPembuat}).\n//\n// We do this the same way a human editor would: click right before the\n// closing \"}\" of each field and type \"Pembuat\".\n\nconst body = context.document.body;\n\n// 1) Insert \"Pembuat\" right after \"{tempatLahir\" (before its closing \"}\").\nlet found = body.search(\"{tempatLahir\", { matchCase: true, matchWholeWord: false });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length === 0) {\n  throw new Error('Could not find \"{tempatLahir\" in the document body.');\n}\nfound.items[0].insertText(\"Pembuat\", Word.InsertLocation.after);\nawait context.sync();\n\n// 2) Insert \"Pembuat\" right after \"}, {ttl\" (before the final closing \"}\").\nfound = body.search(\"}, {ttl\", { matchCase: true, matchWholeWord: false });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length === 0) {\n  throw new Error('Could not find \"}, {ttl\" in the document body.');\n}\nfound.items[0].insertText(\"Pembuat\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# The document contains a merge-field placeholder \"{tempatLahir}, {ttl}\"\n# that must be fixed to \"{tempatLahirPembuat}, {ttlPembuat}\" so the field\n# names match the \"...Pembuat\" convention used elsewhere in the document\n# (e.g. {namaPembuat}, {alamatPembuat}, {pekerjaanPembuat}).\n#\n# We do this the same way a human editor would: click right before the\n# closing \"}\" of each field and type \"Pembuat\".\n\n$d = $word.ActiveDocument\n\n# 1) Insert \"Pembuat\" right after \"{tempatLahir\" (before its closing \"}\").\n$r1 = $d.Content\n$r1.Find.ClearFormatting()\n$r1.Find.Text = \"{tempatLahir\"\n$r1.Find.MatchCase = $true\n$r1.Find.MatchWildcards = $false\n$found1 = $r1.Find.Execute()\nif (-not $found1) {\n    throw 'Could not find \"{tempatLahir\" in the document.'\n}\n$r1.Collapse(0)  # wdCollapseEnd\n$r1.InsertAfter(\"Pembuat\")\n\n# 2) Insert \"Pembuat\" right after \"}, {ttl\" (before the final closing \"}\").\n$r2 = $d.Content\n$r2.Find.ClearFormatting()\n$r2.Find.Text = \"}, {ttl\"\n$r2.Find.MatchCase = $true\n$r2.Find.MatchWildcards = $false\n$found2 = $r2.Find.Execute()\nif (-not $found2) {\n    throw 'Could not find \"}, {ttl\" in the document.'\n}\n$r2.Collapse(0)  # wdCollapseEnd\n$r2.InsertAfter(\"Pembuat\")\n"}
